{"js": "// Word Judgment Entry edit:\n//  1. \"February 01, 2022\" -> \"February 02, 2022\" (both occurrences: the\n//     arraignment date sentence and the \"pay ... in full by\" sentence).\n//  2. \"Defendant waived right to counsel.\" -> \"Defendant was represented\n//     by Dan Smith, Private Counsel.\"\n//  3. The Plea table cell value \"No Contest\" -> \"Guilty\".\n\nconst body = context.document.body;\n\n// 1) Replace every occurrence of the arraignment/payment date.\nconst dateResults = body.search(\"February 01, 2022\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"February 02, 2022\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Counsel sentence.\nconst counselResults = body.search(\"Defendant waived right to counsel.\", { matchCase: true, matchWholeWord: false });\ncounselResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < counselResults.items.length; i++) {\n  counselResults.items[i].insertText(\n    \"Defendant was represented by Dan Smith, Private Counsel.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 3) Plea value in the sentencing table.\nconst pleaResults = body.search(\"No Contest\", { matchCase: true, matchWholeWord: false });\npleaResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < pleaResults.items.length; i++) {\n  pleaResults.items[i].insertText(\"Guilty\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word Judgment Entry edit:\n#  1. \"February 01, 2022\" -> \"February 02, 2022\" (both occurrences: the\n#     arraignment date sentence and the \"pay ... in full by\" sentence).\n#  2. \"Defendant waived right to counsel.\" -> \"Defendant was represented\n#     by Dan Smith, Private Counsel.\"\n#  3. The Plea table cell value \"No Contest\" -> \"Guilty\".\n\n$d = $word.ActiveDocument\n\n# 1) Replace every occurrence of the arraignment/payment date.\n$dateRange = $d.Content\n$dateRange.Find.Execute(\"February 01, 2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"February 02, 2022\", 2)\n\n# 2) Counsel sentence.\n$counselRange = $d.Content\n$counselRange.Find.Execute(\"Defendant waived right to counsel.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Defendant was represented by Dan Smith, Private Counsel.\", 2)\n\n# 3) Plea value in the sentencing table.\n$pleaRange = $d.Content\n$pleaRange.Find.Execute(\"No Contest\", $false, $false, $false, $false, $false, $true, 1, $false, \"Guilty\", 2)\n"}
